# Weekly update: insert a new "Piña" price record for Vega Modelo de Temuco
# as the newest row (row 355), pushing all existing history rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything from row 355 downward by inserting a fresh blank row.
$ws.Rows("355:355").Insert()

# Populate the new row with this week's record.
$ws.Range("A355").Value = 10
$ws.Range("B355").Value = "Vega Modelo de Temuco"
$ws.Range("C355").Value = "La Araucanía"
$ws.Range("D355").Value = 44585
$ws.Range("E355").Value = 9
$ws.Range("F355").Value = "Fruta"
$ws.Range("G355").Value = 100108
$ws.Range("H355").Value = "Tropicales y subtropicales"
$ws.Range("I355").Value = 100108005
$ws.Range("J355").Value = "Piña"
$ws.Range("K355").Value = "Caramelo"
$ws.Range("L355").Value = "Primera"
$ws.Range("M355").Value = 200
$ws.Range("N355").Value = 18000
$ws.Range("O355").Value = 18000
$ws.Range("P355").Value = 18000
$ws.Range("Q355").Value = "$/caja 12 unidades"
$ws.Range("R355").Value = "Ecuador"
$ws.Range("S355").Value = 1500
$ws.Range("T355").Value = 12
